# Auto-generated Excel COM-interop script to apply the workbook diff
# Updates numeric cells (currentAveragePrice* / LevePrice* / LeveProfit* columns)
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3413.3635
$ws.Range("I69").Value = 3381.5
$ws.Range("J69").Value = 3498.3333
$ws.Range("K69").Value = 10144.5
$ws.Range("L69").Value = 10494.9999
$ws.Range("M69").Value = -9270.5
$ws.Range("N69").Value = -12242.9999
$ws.Range("H72").Value = 3413.3635
$ws.Range("I72").Value = 3381.5
$ws.Range("J72").Value = 3498.3333
$ws.Range("K72").Value = 30433.5
$ws.Range("L72").Value = 31484.9997
$ws.Range("M72").Value = -26065.5
$ws.Range("N72").Value = -40220.9997
$ws.Range("H107").Value = 7978.857
$ws.Range("I107").Value = 8584.462
$ws.Range("J107").Value = 106
$ws.Range("K107").Value = 8584.462
$ws.Range("L107").Value = 106
$ws.Range("M107").Value = -6664.462
$ws.Range("N107").Value = -3946
$ws.Range("H112").Value = 1618.8684
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 1809.9062
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 5429.7186
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -7645.7186
$ws.Range("H116").Value = 3027.2727
$ws.Range("I116").Value = 2442.8572
$ws.Range("K116").Value = 2442.8572
$ws.Range("M116").Value = 999.1428000000001
$ws.Range("H129").Value = 1535.6538
$ws.Range("I129").Value = 489.33334
$ws.Range("J129").Value = 1754.6511
$ws.Range("K129").Value = 1468.00002
$ws.Range("L129").Value = 5263.9533
$ws.Range("M129").Value = 3531.99998
$ws.Range("N129").Value = -15263.9533

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11022.032
$ws.Range("I32").Value = 10166.988
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 10166.988
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -9879.987999999999
$ws.Range("N32").Value = -20574
$ws.Range("H34").Value = 15305.714
$ws.Range("I34").Value = 11000
$ws.Range("J34").Value = 16023.333
$ws.Range("K34").Value = 11000
$ws.Range("L34").Value = 16023.333
$ws.Range("M34").Value = -10729
$ws.Range("N34").Value = -16565.333
$ws.Range("H74").Value = 1841.079
$ws.Range("I74").Value = 1131.7407
$ws.Range("J74").Value = 3582.182
$ws.Range("K74").Value = 1131.7407
$ws.Range("L74").Value = 3582.182
$ws.Range("M74").Value = -257.7407000000001
$ws.Range("N74").Value = -5330.182
$ws.Range("H77").Value = 1841.079
$ws.Range("I77").Value = 1131.7407
$ws.Range("J77").Value = 3582.182
$ws.Range("K77").Value = 5658.703500000001
$ws.Range("L77").Value = 17910.91
$ws.Range("M77").Value = -1290.703500000001
$ws.Range("N77").Value = -26646.91
$ws.Range("H132").Value = 5001997
$ws.Range("I132").Value = 6946311
$ws.Range("J132").Value = 2332.3572
$ws.Range("K132").Value = 20838933
$ws.Range("L132").Value = 6997.071599999999
$ws.Range("M132").Value = -20836403
$ws.Range("N132").Value = -12057.0716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 31252870
$ws.Range("I86").Value = 2575.5
$ws.Range("J86").Value = 125003750
$ws.Range("K86").Value = 2575.5
$ws.Range("L86").Value = 125003750
$ws.Range("M86").Value = -1452.5
$ws.Range("N86").Value = -125005996
$ws.Range("H89").Value = 31252870
$ws.Range("I89").Value = 2575.5
$ws.Range("J89").Value = 125003750
$ws.Range("K89").Value = 12877.5
$ws.Range("L89").Value = 625018750
$ws.Range("M89").Value = -7261.5
$ws.Range("N89").Value = -625029982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2199.131
$ws.Range("I31").Value = 2103.3635
$ws.Range("J31").Value = 2253.1538
$ws.Range("K31").Value = 2103.3635
$ws.Range("L31").Value = 2253.1538
$ws.Range("M31").Value = -1808.3635
$ws.Range("N31").Value = -2843.1538
$ws.Range("H34").Value = 2199.131
$ws.Range("I34").Value = 2103.3635
$ws.Range("J34").Value = 2253.1538
$ws.Range("K34").Value = 2103.3635
$ws.Range("L34").Value = 2253.1538
$ws.Range("M34").Value = -1901.3635
$ws.Range("N34").Value = -2657.1538
$ws.Range("H99").Value = 1830.3846
$ws.Range("I99").Value = 1830.3846
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1830.3846
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -332.3846000000001
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 1011.3939
$ws.Range("I105").Value = 909.25
$ws.Range("J105").Value = 1283.7778
$ws.Range("K105").Value = 909.25
$ws.Range("L105").Value = 1283.7778
$ws.Range("M105").Value = 837.75
$ws.Range("N105").Value = -4777.7778
$ws.Range("H126").Value = 1830.3846
$ws.Range("I126").Value = 1830.3846
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5491.1538
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3021.1538
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 440.51852
$ws.Range("I107").Value = 269.125
$ws.Range("J107").Value = 512.6842
$ws.Range("K107").Value = 807.375
$ws.Range("L107").Value = 1538.0526
$ws.Range("M107").Value = 1112.625
$ws.Range("N107").Value = -5378.0526
$ws.Range("H118").Value = 1746.862
$ws.Range("I118").Value = 701.4
$ws.Range("J118").Value = 1964.6666
$ws.Range("K118").Value = 2104.2
$ws.Range("L118").Value = 5893.9998
$ws.Range("M118").Value = -861.1999999999998
$ws.Range("N118").Value = -8379.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1480.1177
$ws.Range("I97").Value = 1303.25
$ws.Range("J97").Value = 1904.6
$ws.Range("K97").Value = 1303.25
$ws.Range("L97").Value = 1904.6
$ws.Range("M97").Value = -807.25
$ws.Range("N97").Value = -2896.6
$ws.Range("H102").Value = 4595.75
$ws.Range("I102").Value = 5730.364
$ws.Range("K102").Value = 5730.364
$ws.Range("M102").Value = -4108.364
$ws.Range("H126").Value = 3989.5881
$ws.Range("I126").Value = 2202
$ws.Range("J126").Value = 4964.636
$ws.Range("K126").Value = 6606
$ws.Range("L126").Value = 14893.908
$ws.Range("M126").Value = -4136
$ws.Range("N126").Value = -19833.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5481.129
$ws.Range("I7").Value = 5909.643
$ws.Range("J7").Value = 5128.2354
$ws.Range("K7").Value = 5909.643
$ws.Range("L7").Value = 5128.2354
$ws.Range("M7").Value = -5797.643
$ws.Range("N7").Value = -5352.2354
$ws.Range("H43").Value = 19300
$ws.Range("J43").Value = 19300
$ws.Range("L43").Value = 19300
$ws.Range("N43").Value = -19686
$ws.Range("H122").Value = 9389.941000000001
$ws.Range("I122").Value = 10137.667
$ws.Range("K122").Value = 30413.001
$ws.Range("M122").Value = -27963.001
$ws.Range("H126").Value = 5481.129
$ws.Range("I126").Value = 5909.643
$ws.Range("J126").Value = 5128.2354
$ws.Range("K126").Value = 17728.929
$ws.Range("L126").Value = 15384.7062
$ws.Range("M126").Value = -15258.929
$ws.Range("N126").Value = -20324.7062
$ws.Range("H132").Value = 13521845
$ws.Range("I132").Value = 5224.8213
$ws.Range("J132").Value = 55573550
$ws.Range("K132").Value = 15674.4639
$ws.Range("L132").Value = 166720650
$ws.Range("M132").Value = -13144.4639
$ws.Range("N132").Value = -166725710
$ws.Range("H136").Value = 4588.421
$ws.Range("I136").Value = 2058.7144
$ws.Range("J136").Value = 34101.668
$ws.Range("K136").Value = 6176.1432
$ws.Range("L136").Value = 102305.004
$ws.Range("M136").Value = -3626.1432
$ws.Range("N136").Value = -107405.004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 14524.692
$ws.Range("I33").Value = 5500
$ws.Range("J33").Value = 16165.546
$ws.Range("K33").Value = 5500
$ws.Range("L33").Value = 16165.546
$ws.Range("M33").Value = -5250
$ws.Range("N33").Value = -16665.546
$ws.Range("H36").Value = 14524.692
$ws.Range("I36").Value = 5500
$ws.Range("J36").Value = 16165.546
$ws.Range("K36").Value = 5500
$ws.Range("L36").Value = 16165.546
$ws.Range("M36").Value = -5250
$ws.Range("N36").Value = -16665.546
$ws.Range("H37").Value = 17645
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 17645
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 17645
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -18051
$ws.Range("H41").Value = 7550.6
$ws.Range("I41").Value = 4749.5
$ws.Range("J41").Value = 9418
$ws.Range("K41").Value = 4749.5
$ws.Range("L41").Value = 9418
$ws.Range("M41").Value = -4359.5
$ws.Range("N41").Value = -10198
$ws.Range("H126").Value = 1684.6818
$ws.Range("I126").Value = 1891.9412
$ws.Range("J126").Value = 980
$ws.Range("K126").Value = 5675.8236
$ws.Range("L126").Value = 2940
$ws.Range("M126").Value = -3205.8236
$ws.Range("N126").Value = -7880
